$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(100).Insert()

$ws.Range("A100").Value = 2
$ws.Range("B100").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C100").Value = "Coquimbo"
$ws.Range("D100").Value = 44994
$ws.Range("E100").Value = 4
$ws.Range("F100").Value = 100112043
$ws.Range("G100").Value = "Pepino ensalada"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = 6500
$ws.Range("N100").Value = "`$/caja 70 unidades"
$ws.Range("O100").Value = "Provincia de Limarí"
$ws.Range("P100").Value = 93
$ws.Range("Q100").Value = 70
$ws.Range("R100").Value = "Hortaliza"
